# TC13_CDS_Filter_PHSAccession-phs003111.xlsx edit
# - Fix a dropped comma in the ParticipantsTab Cypher query (B2)
# - Bump the "query" column font from 12pt (theme/automatic color) to a
#   14pt explicit-black font across the query + StatQuery columns (B2:C4)
# - Move the active selection from B2 to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the ParticipantsTab query text in B2 -----------------------
# The "samples" field used to coalesce the joined text with a trailing
# ", \"\"" arg; the fixed query drops that one comma.
$old = 'samples: coalesce(apoc.text.join(apoc.coll.sort(collect(distinct samp.sample_id)), ", "), "")'
$new = 'samples: coalesce(apoc.text.join(apoc.coll.sort(collect(distinct samp.sample_id)), ", ") "")'

$participantsQuery = $ws.Range("B2").Value2
if ($participantsQuery.Contains($old)) {
    $ws.Range("B2").Value = $participantsQuery.Replace($old, $new)
}

# Re-writing the (still enormous) wrapped string nudges the recalculated
# wrap height away from Excel's real 409.6pt row-height cap; pin it back so
# row 2 keeps matching rows 3 & 4.
$ws.Rows.Item(2).RowHeight = 409.6

# --- 2. Re-style the query / StatQuery columns ------------------------------
# Both columns share the same bumped font: 14pt, explicit black (was 12pt
# with the automatic/theme text color).
$queryCells = $ws.Range("B2:C4")
$queryCells.Font.Size = 14
$queryCells.Font.Color = 0

# --- 3. Move the selection to D2 --------------------------------------------
$ws.Range("D2").Select()
